$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.917.47"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.550.39"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").Value = "206.47"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "0.489"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("D8").Value = "21.97"
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").Value = "0.0595"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").Value = "1.770.23"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "1.534.18"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").Value = "3.75"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "26.888.55"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "61.59"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "0.0₃0712"
$ws.Range("E18").Value = "  +3.40%  "
$ws.Range("D19").Value = "216.71"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").Value = "7.28"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").Value = "9.19"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "1.95"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").Value = "153.26"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "14.95"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("D34").Value = "1.404.30"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("E35").Value = "  +1.98%  "
$ws.Range("D36").Value = "0.959"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("D39").Value = "0.525"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").Value = "5.65"
$ws.Range("E42").Value = "  +3.89%  "
$ws.Range("E43").Value = "  +1.28%  "
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("D45").Value = "64.40"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("D47").Value = "1.684.15"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("D48").Value = "87.28"
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("E50").Value = "  +5.81%  "
$ws.Range("E51").Value = "  -0.06%  "
